{"js": "// Update the three-digit-by-one-digit multiplication prompts in the\n// practice table. Each \"old\" expression is unique in the document, so a\n// simple search-and-replace per pair is safe and order independent.\nconst replacements = [\n  [\"650\u00d75=\", \"629\u00d77=\"],\n  [\"279\u00d78=\", \"427\u00d76=\"],\n  [\"566\u00d78=\", \"678\u00d77=\"],\n  [\"756\u00d72=\", \"679\u00d76=\"],\n  [\"664\u00d75=\", \"394\u00d72=\"],\n  [\"444\u00d76=\", \"339\u00d78=\"],\n  [\"898\u00d74=\", \"272\u00d79=\"],\n  [\"940\u00d77=\", \"919\u00d77=\"],\n  [\"297\u00d74=\", \"685\u00d75=\"],\n  [\"258\u00d79=\", \"729\u00d78=\"],\n  [\"118\u00d73=\", \"142\u00d76=\"],\n  [\"788\u00d79=\", \"288\u00d78=\"],\n  [\"458\u00d75=\", \"746\u00d78=\"],\n  [\"898\u00d78=\", \"300\u00d72=\"],\n  [\"583\u00d79=\", \"331\u00d78=\"],\n  [\"834\u00d73=\", \"443\u00d73=\"],\n  [\"577\u00d77=\", \"770\u00d74=\"],\n  [\"606\u00d78=\", \"612\u00d79=\"],\n  [\"854\u00d74=\", \"890\u00d74=\"],\n  [\"523\u00d78=\", \"796\u00d79=\"],\n  [\"323\u00d76=\", \"778\u00d74=\"],\n  [\"822\u00d74=\", \"396\u00d74=\"],\n  [\"930\u00d72=\", \"624\u00d75=\"],\n  [\"387\u00d75=\", \"999\u00d75=\"],\n  [\"842\u00d73=\", \"537\u00d75=\"],\n];\n\nconst body = context.document.body;\n\n// Kick off a search for every old value up front, then sync once so all\n// results are loaded together.\nconst searches = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearches.forEach((s) => s.load(\"items\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searches[i].items;\n  for (let j = 0; j < items.length; j++) {\n    items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit-by-one-digit multiplication prompts in the\n# practice table. Each \"old\" expression occurs exactly once in the\n# document, so a plain Find/Replace per pair (no \"replace all\" ambiguity)\n# is sufficient and keeps the original run formatting untouched.\n$pairs = @(\n    @{ Old = \"650\u00d75=\"; New = \"629\u00d77=\" },\n    @{ Old = \"279\u00d78=\"; New = \"427\u00d76=\" },\n    @{ Old = \"566\u00d78=\"; New = \"678\u00d77=\" },\n    @{ Old = \"756\u00d72=\"; New = \"679\u00d76=\" },\n    @{ Old = \"664\u00d75=\"; New = \"394\u00d72=\" },\n    @{ Old = \"444\u00d76=\"; New = \"339\u00d78=\" },\n    @{ Old = \"898\u00d74=\"; New = \"272\u00d79=\" },\n    @{ Old = \"940\u00d77=\"; New = \"919\u00d77=\" },\n    @{ Old = \"297\u00d74=\"; New = \"685\u00d75=\" },\n    @{ Old = \"258\u00d79=\"; New = \"729\u00d78=\" },\n    @{ Old = \"118\u00d73=\"; New = \"142\u00d76=\" },\n    @{ Old = \"788\u00d79=\"; New = \"288\u00d78=\" },\n    @{ Old = \"458\u00d75=\"; New = \"746\u00d78=\" },\n    @{ Old = \"898\u00d78=\"; New = \"300\u00d72=\" },\n    @{ Old = \"583\u00d79=\"; New = \"331\u00d78=\" },\n    @{ Old = \"834\u00d73=\"; New = \"443\u00d73=\" },\n    @{ Old = \"577\u00d77=\"; New = \"770\u00d74=\" },\n    @{ Old = \"606\u00d78=\"; New = \"612\u00d79=\" },\n    @{ Old = \"854\u00d74=\"; New = \"890\u00d74=\" },\n    @{ Old = \"523\u00d78=\"; New = \"796\u00d79=\" },\n    @{ Old = \"323\u00d76=\"; New = \"778\u00d74=\" },\n    @{ Old = \"822\u00d74=\"; New = \"396\u00d74=\" },\n    @{ Old = \"930\u00d72=\"; New = \"624\u00d75=\" },\n    @{ Old = \"387\u00d75=\"; New = \"999\u00d75=\" },\n    @{ Old = \"842\u00d73=\"; New = \"537\u00d75=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceOne = 1\n    $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 1) | Out-Null\n}\n"}
